$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 524.875
$ws.Range("I12").Value = 581.4286
$ws.Range("K12").Value = 581.4286
$ws.Range("M12").Value = -411.4286
$ws.Range("H19").Value = 1461.3529
$ws.Range("I19").Value = 1457.4286
$ws.Range("K19").Value = 1457.4286
$ws.Range("M19").Value = -1282.4286
$ws.Range("H21").Value = 119998
$ws.Range("I21").Value = 119998
$ws.Range("K21").Value = 119998
$ws.Range("M21").Value = -119530
$ws.Range("H23").Value = 119998
$ws.Range("I23").Value = 119998
$ws.Range("K23").Value = 119998
$ws.Range("M23").Value = -119764
$ws.Range("H28").Value = 1167.3043
$ws.Range("I28").Value = 1128.8422
$ws.Range("J28").Value = 1350
$ws.Range("K28").Value = 1128.8422
$ws.Range("L28").Value = 1350
$ws.Range("M28").Value = -643.8422
$ws.Range("N28").Value = -2320
$ws.Range("H38").Value = 3369.3333
$ws.Range("I38").Value = 3369.3333
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 10107.9999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -9735.999899999999
$ws.Range("N38").ClearContents()
$ws.Range("H51").Value = 3797.0715
$ws.Range("I51").Value = 4151.6313
$ws.Range("J51").Value = 3048.5557
$ws.Range("K51").Value = 4151.6313
$ws.Range("L51").Value = 3048.5557
$ws.Range("M51").Value = -3667.6313
$ws.Range("N51").Value = -4016.5557
$ws.Range("H58").Value = 240.66667
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H87").Value = 110085.25
$ws.Range("J87").Value = 91784.664
$ws.Range("L87").Value = 91784.664
$ws.Range("N87").Value = -94280.664
$ws.Range("H90").Value = 110085.25
$ws.Range("J90").Value = 91784.664
$ws.Range("L90").Value = 275353.992
$ws.Range("N90").Value = -287833.992
$ws.Range("H107").Value = 1658
$ws.Range("J107").Value = 1600
$ws.Range("L107").Value = 1600
$ws.Range("N107").Value = -5440
$ws.Range("H116").Value = 12023.5625
$ws.Range("I116").Value = 20749.125
$ws.Range("J116").Value = 3298
$ws.Range("K116").Value = 20749.125
$ws.Range("L116").Value = 3298
$ws.Range("M116").Value = -17307.125
$ws.Range("N116").Value = -10182
$ws.Range("H132").Value = 1636.0204
$ws.Range("I132").Value = 1348.6571
$ws.Range("J132").Value = 2354.4285
$ws.Range("K132").Value = 4045.9713
$ws.Range("L132").Value = 7063.2855
$ws.Range("M132").Value = -1515.9713
$ws.Range("N132").Value = -12123.2855
$ws.Range("H141").Value = 5641.0938
$ws.Range("I141").Value = 4833.8887
$ws.Range("K141").Value = 14501.6661
$ws.Range("M141").Value = -9321.666100000002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1746.6923
$ws.Range("I74").Value = 1630.1936
$ws.Range("K74").Value = 1630.1936
$ws.Range("M74").Value = -756.1936000000001
$ws.Range("H77").Value = 1746.6923
$ws.Range("I77").Value = 1630.1936
$ws.Range("K77").Value = 8150.968000000001
$ws.Range("M77").Value = -3782.968000000001
$ws.Range("H88").Value = 1066.8334
$ws.Range("I88").Value = 945.7143
$ws.Range("K88").Value = 945.7143
$ws.Range("M88").Value = -539.7143
$ws.Range("H91").Value = 1066.8334
$ws.Range("I91").Value = 945.7143
$ws.Range("K91").Value = 945.7143
$ws.Range("M91").Value = 458.2857
$ws.Range("H109").Value = 54770
$ws.Range("J109").Value = 54770
$ws.Range("L109").Value = 54770
$ws.Range("N109").Value = -57544
$ws.Range("H110").Value = 3008.8
$ws.Range("I110").Value = 2955.5715
$ws.Range("J110").Value = 3133
$ws.Range("K110").Value = 2955.5715
$ws.Range("L110").Value = 3133
$ws.Range("M110").Value = -910.5715
$ws.Range("N110").Value = -7223
$ws.Range("H132").Value = 2767.457
$ws.Range("I132").Value = 2719.724
$ws.Range("K132").Value = 8159.172
$ws.Range("M132").Value = -5629.172
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H80").Value = 590.6
$ws.Range("J80").Value = 686.2857
$ws.Range("L80").Value = 686.2857
$ws.Range("N80").Value = -2682.2857
$ws.Range("H83").Value = 590.6
$ws.Range("J83").Value = 686.2857
$ws.Range("L83").Value = 3431.4285
$ws.Range("N83").Value = -13415.4285
$ws.Range("H86").Value = 755.1177
$ws.Range("I86").Value = 750
$ws.Range("J86").Value = 764.5
$ws.Range("K86").Value = 750
$ws.Range("L86").Value = 764.5
$ws.Range("M86").Value = 373
$ws.Range("N86").Value = -3010.5
$ws.Range("H89").Value = 755.1177
$ws.Range("I89").Value = 750
$ws.Range("J89").Value = 764.5
$ws.Range("K89").Value = 3750
$ws.Range("L89").Value = 3822.5
$ws.Range("M89").Value = 1866
$ws.Range("N89").Value = -15054.5
$ws.Range("H94").Value = 582.3929000000001
$ws.Range("I94").Value = 437.9
$ws.Range("K94").Value = 437.9
$ws.Range("M94").Value = 13.10000000000002
$ws.Range("H116").Value = 99749.5
$ws.Range("J116").Value = 99749.5
$ws.Range("L116").Value = 99749.5
$ws.Range("N116").Value = -108927.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 20.76923
$ws.Range("I7").Value = 15.24
$ws.Range("K7").Value = 15.24
$ws.Range("M7").Value = 97.76000000000001
$ws.Range("H31").Value = 4397.533
$ws.Range("I31").Value = 1906.7142
$ws.Range("K31").Value = 1906.7142
$ws.Range("M31").Value = -1611.7142
$ws.Range("H34").Value = 4397.533
$ws.Range("I34").Value = 1906.7142
$ws.Range("K34").Value = 1906.7142
$ws.Range("M34").Value = -1704.7142
$ws.Range("H58").Value = 2620.5789
$ws.Range("I58").Value = 2430.82
$ws.Range("J58").Value = 3976
$ws.Range("K58").Value = 2430.82
$ws.Range("L58").Value = 3976
$ws.Range("M58").Value = -2227.82
$ws.Range("N58").Value = -4382
$ws.Range("H68").Value = 69993.336
$ws.Range("J68").Value = 69993.336
$ws.Range("L68").Value = 69993.336
$ws.Range("N68").Value = -71491.336
$ws.Range("H71").Value = 69993.336
$ws.Range("J71").Value = 69993.336
$ws.Range("L71").Value = 209980.008
$ws.Range("N71").Value = -217468.008
$ws.Range("H74").Value = 59999.668
$ws.Range("J74").Value = 59999.668
$ws.Range("L74").Value = 59999.668
$ws.Range("N74").Value = -61747.668
$ws.Range("H77").Value = 59999.668
$ws.Range("J77").Value = 59999.668
$ws.Range("L77").Value = 179999.004
$ws.Range("N77").Value = -188735.004
$ws.Range("H98").Value = 63369.5
$ws.Range("J98").Value = 63369.5
$ws.Range("L98").Value = 63369.5
$ws.Range("N98").Value = -67861.5
$ws.Range("H99").Value = 2809.875
$ws.Range("I99").Value = 2789
$ws.Range("K99").Value = 2789
$ws.Range("M99").Value = -1291
$ws.Range("H112").Value = 73381.836
$ws.Range("J112").Value = 73381.836
$ws.Range("L112").Value = 73381.836
$ws.Range("N112").Value = -76335.836
$ws.Range("H119").Value = 116998.5
$ws.Range("J119").Value = 116998.5
$ws.Range("L119").Value = 116998.5
$ws.Range("N119").Value = -126674.5
$ws.Range("H122").Value = 5251.8667
$ws.Range("I122").Value = 4174.125
$ws.Range("J122").Value = 6483.5713
$ws.Range("K122").Value = 12522.375
$ws.Range("L122").Value = 19450.7139
$ws.Range("M122").Value = -10072.375
$ws.Range("N122").Value = -24350.7139
$ws.Range("H126").Value = 2809.875
$ws.Range("I126").Value = 2789
$ws.Range("K126").Value = 8367
$ws.Range("M126").Value = -5897
$ws.Range("H132").Value = 3781.9583
$ws.Range("I132").Value = 3026.7368
$ws.Range("J132").Value = 6651.8
$ws.Range("K132").Value = 9080.2104
$ws.Range("L132").Value = 19955.4
$ws.Range("M132").Value = -6550.2104
$ws.Range("N132").Value = -25015.4
$ws.Range("H134").Value = 1865
$ws.Range("I134").Value = 1624.6
$ws.Range("K134").Value = 4873.799999999999
$ws.Range("M134").Value = -2338.799999999999
$ws.Range("H136").Value = 2620.5789
$ws.Range("I136").Value = 2430.82
$ws.Range("J136").Value = 3976
$ws.Range("K136").Value = 7292.460000000001
$ws.Range("L136").Value = 11928
$ws.Range("M136").Value = -4742.460000000001
$ws.Range("N136").Value = -17028
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 865.2222
$ws.Range("I68").Value = 801.5
$ws.Range("J68").Value = 883.4286
$ws.Range("K68").Value = 2404.5
$ws.Range("L68").Value = 2650.2858
$ws.Range("M68").Value = -1593.5
$ws.Range("N68").Value = -4272.2858
$ws.Range("H71").Value = 865.2222
$ws.Range("I71").Value = 801.5
$ws.Range("J71").Value = 883.4286
$ws.Range("K71").Value = 7213.5
$ws.Range("L71").Value = 7950.8574
$ws.Range("M71").Value = -3157.5
$ws.Range("N71").Value = -16062.8574
$ws.Range("H107").Value = 1120.5333
$ws.Range("I107").Value = 1167.7142
$ws.Range("K107").Value = 3503.1426
$ws.Range("M107").Value = -1583.1426
$ws.Range("H120").Value = 21845.924
$ws.Range("I120").Value = 11332.333
$ws.Range("K120").Value = 33996.999
$ws.Range("M120").Value = -29158.999
$ws.Range("H121").Value = 11548.9
$ws.Range("J121").Value = 18788.334
$ws.Range("L121").Value = 56365.00199999999
$ws.Range("N121").Value = -58985.00199999999
$ws.Range("H122").Value = 965.6667
$ws.Range("I122").Value = 949.5
$ws.Range("K122").Value = 8545.5
$ws.Range("M122").Value = -6095.5
$ws.Range("H129").Value = 1938.9
$ws.Range("J129").Value = 3027.3333
$ws.Range("L129").Value = 9081.999899999999
$ws.Range("N129").Value = -19081.9999
$ws.Range("H138").Value = 35501148
$ws.Range("I138").Value = 1722.5
$ws.Range("K138").Value = 5167.5
$ws.Range("M138").Value = -27.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3151.1667
$ws.Range("I126").Value = 2870.5
$ws.Range("K126").Value = 8611.5
$ws.Range("M126").Value = -6141.5
$ws.Range("H132").Value = 3134.9333
$ws.Range("I132").Value = 2786.1304
$ws.Range("J132").Value = 4281
$ws.Range("K132").Value = 8358.3912
$ws.Range("L132").Value = 12843
$ws.Range("M132").Value = -5828.3912
$ws.Range("N132").Value = -17903
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2774.8
$ws.Range("J40").Value = 2974.8333
$ws.Range("L40").Value = 2974.8333
$ws.Range("N40").Value = -3246.8333
$ws.Range("H64").Value = 10000
$ws.Range("I64").Value = 10000
$ws.Range("K64").Value = 10000
$ws.Range("M64").Value = -9775
$ws.Range("H67").Value = 10000
$ws.Range("I67").Value = 10000
$ws.Range("K67").Value = 10000
$ws.Range("M67").Value = -9220
$ws.Range("H88").Value = 76021.125
$ws.Range("I88").Value = 69595.57000000001
$ws.Range("J88").Value = 121000
$ws.Range("K88").Value = 69595.57000000001
$ws.Range("L88").Value = 121000
$ws.Range("M88").Value = -69167.57000000001
$ws.Range("N88").Value = -121856
$ws.Range("H91").Value = 76021.125
$ws.Range("I91").Value = 69595.57000000001
$ws.Range("J91").Value = 121000
$ws.Range("K91").Value = 69595.57000000001
$ws.Range("L91").Value = 121000
$ws.Range("M91").Value = -68113.57000000001
$ws.Range("N91").Value = -123964
$ws.Range("H132").Value = 4937
$ws.Range("I132").Value = 3895.6667
$ws.Range("K132").Value = 11687.0001
$ws.Range("M132").Value = -9157.000100000001
$ws.Range("H136").Value = 5073.1035
$ws.Range("I136").Value = 3322.611
$ws.Range("J136").Value = 7937.5454
$ws.Range("K136").Value = 9967.832999999999
$ws.Range("L136").Value = 23812.6362
$ws.Range("M136").Value = -7417.832999999999
$ws.Range("N136").Value = -28912.6362
$ws.Range("H139").Value = 100870.836
$ws.Range("J139").Value = 106045.4
$ws.Range("L139").Value = 106045.4
$ws.Range("N139").Value = -116325.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 13499
$ws.Range("I8").Value = 13499
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 13499
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -13359
$ws.Range("N8").ClearContents()
$ws.Range("H59").Value = 49999
$ws.Range("J59").Value = 49999
$ws.Range("L59").Value = 49999
$ws.Range("N59").Value = -51475
$ws.Range("H62").Value = 6068.625
$ws.Range("J62").Value = 7008.1665
$ws.Range("L62").Value = 7008.1665
$ws.Range("N62").Value = -8256.166499999999
$ws.Range("H65").Value = 6068.625
$ws.Range("J65").Value = 7008.1665
$ws.Range("L65").Value = 35040.8325
$ws.Range("N65").Value = -41280.8325
$ws.Range("H81").Value = 3275.3914
$ws.Range("I81").Value = 1811.1538
$ws.Range("K81").Value = 3622.3076
$ws.Range("M81").Value = -2561.3076
$ws.Range("H82").Value = 90000
$ws.Range("J82").Value = 90000
$ws.Range("L82").Value = 90000
$ws.Range("N82").Value = -90766
$ws.Range("H84").Value = 3275.3914
$ws.Range("I84").Value = 1811.1538
$ws.Range("K84").Value = 18111.538
$ws.Range("M84").Value = -12807.538
$ws.Range("H85").Value = 90000
$ws.Range("J85").Value = 90000
$ws.Range("L85").Value = 90000
$ws.Range("N85").Value = -92652
$ws.Range("H96").Value = 11995.846
$ws.Range("I96").Value = 7427.7144
$ws.Range("K96").Value = 7427.7144
$ws.Range("M96").Value = -6054.7144
$ws.Range("H126").Value = 2570
$ws.Range("I126").Value = 966.6667
$ws.Range("K126").Value = 2900.0001
$ws.Range("M126").Value = -430.0001000000002
$ws.Range("H132").Value = 1949.5
$ws.Range("I132").Value = 1708.4
$ws.Range("J132").Value = 3155
$ws.Range("K132").Value = 5125.200000000001
$ws.Range("L132").Value = 9465
$ws.Range("M132").Value = -2595.200000000001
$ws.Range("N132").Value = -14525
